# Updated capital structure database
# Applies the per-cell value updates to the "oman_insurance_general" sheet,
# including a brand-new row 6 for Takaful Oman Insurance SAOG (MSM:TAOI).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# --- Row 2 - Oman Insurance & Reinsurance Co. (placeholder name '4') ---
$c = $ws.Cells.Item(2,2)  # B2: keep as TEXT "4" (not a number)
$c.NumberFormat = "@"
$c.Value = "4"
$c.Style = "Normal"
$ws.Cells.Item(2,4).Value = 0.01755  # D2
$ws.Cells.Item(2,5).Value = 0.174  # E2
$ws.Cells.Item(2,7).Value = 0.1318227593152064  # G2
$ws.Cells.Item(2,8).Value = 0.1318227593152064  # H2
$ws.Cells.Item(2,9).Value = 0.165508559919436  # I2
$ws.Cells.Item(2,10).Value = 0.1456953753912002  # J2
$ws.Cells.Item(2,11).Value = 25.16  # K2
$ws.Cells.Item(2,12).Value = 0.126686807653575  # L2
$ws.Cells.Item(2,13).Value = 10.52  # M2
$ws.Cells.Item(2,14).Value = 0.0492970946579194  # N2
$ws.Cells.Item(2,15).Value = 0.4181240063593005  # O2
$ws.Cells.Item(2,16).Value = 10.52  # P2
$ws.Cells.Item(2,17).Value = 0.0492970946579194  # Q2
$ws.Cells.Item(2,18).Value = 0.4181240063593005  # R2
$ws.Cells.Item(2,21).Value = 35.353  # U2
$ws.Cells.Item(2,22).Value = 0.1656654170571696  # V2
$ws.Cells.Item(2,23).Value = 0.1099820578289973  # W2
$ws.Cells.Item(2,24).Value = 0.06600886964855902  # X2
$ws.Cells.Item(2,25).Value = 0.04397318818043829  # Y2
$ws.Cells.Item(2,26).Value = 0.8519583029471066  # Z2
$ws.Cells.Item(2,27).Value = 0.1277407710708613  # AA2
$ws.Cells.Item(2,28).Value = 0.06523099898918069  # AB2
$ws.Cells.Item(2,29).Value = 0.06415522656919559  # AC2
$ws.Cells.Item(2,30).Value = 16.634  # AD2
$ws.Cells.Item(2,31).Value = 0  # AE2
$ws.Cells.Item(2,32).Value = 16.634  # AF2
$ws.Cells.Item(2,33).Value = -18.719  # AG2
$ws.Cells.Item(2,34).Value = 0.07231104967091821  # AH2
$ws.Cells.Item(2,35).Value = 0.05856341142257617  # AI2
$ws.Cells.Item(2,36).Value = -0.0961521668781237  # AJ2
$ws.Cells.Item(2,37).Value = -0.07527314109240353  # AK2
$ws.Cells.Item(2,38).Value = 1.21  # AL2
$ws.Cells.Item(2,39).Value = 1.21  # AM2
$ws.Cells.Item(2,40).Value = 0.4775768016078094  # AN2
$ws.Cells.Item(2,41).Value = 27.16528925619835  # AO2
$ws.Cells.Item(2,42).Value = -0.537438989376974  # AP2
$ws.Cells.Item(2,43).Value = 27.16528925619835  # AQ2

# --- Row 3 - renamed to Dhofar Insurance Company SAOG (MSM:DICS) ---
$ws.Cells.Item(3,2).Value = "Dhofar Insurance Company SAOG (MSM:DICS)"  # B3
$ws.Cells.Item(3,4).Value = -0.156  # D3
$ws.Cells.Item(3,5).Value = 0.146  # E3
$ws.Cells.Item(3,7).Value = 0.1676557863501484  # G3
$ws.Cells.Item(3,8).Value = 0.1676557863501484  # H3
$ws.Cells.Item(3,9).Value = 0.2077151335311573  # I3
$ws.Cells.Item(3,10).Value = 0.1787493534424087  # J3
$ws.Cells.Item(3,11).Value = 9.359999999999999  # K3
$ws.Cells.Item(3,12).Value = 0.1388724035608308  # L3
$ws.Cells.Item(3,21).Value = 17.3  # U3
$ws.Cells.Item(3,22).Value = 0.4061032863849765  # V3
$ws.Cells.Item(3,23).Value = 0.1674418604651163  # W3
$ws.Cells.Item(3,24).Value = 0.077606546158449  # X3
$ws.Cells.Item(3,25).Value = 0.08983531430666727  # Y3
$ws.Cells.Item(3,26).Value = 1.227686703096539  # Z3
$ws.Cells.Item(3,27).Value = 0.2194482044083488  # AA3
$ws.Cells.Item(3,28).Value = 0.06959898548324508  # AB3
$ws.Cells.Item(3,29).Value = 0.1498492189251037  # AC3
$ws.Cells.Item(3,30).Value = 13  # AD3
$ws.Cells.Item(3,31).Value = 0  # AE3
$ws.Cells.Item(3,32).Value = 13  # AF3
$ws.Cells.Item(3,33).Value = -4.300000000000001  # AG3
$ws.Cells.Item(3,34).Value = 0.2338129496402878  # AH3
$ws.Cells.Item(3,35).Value = 0.1681759379042691  # AI3
$ws.Cells.Item(3,36).Value = -0.1122715404699739  # AJ3
$ws.Cells.Item(3,37).Value = -0.07166666666666668  # AK3
$ws.Cells.Item(3,38).Value = 1.21  # AL3
$ws.Cells.Item(3,39).Value = 1.21  # AM3
$ws.Cells.Item(3,40).Value = 0.9027777777777778  # AN3
$ws.Cells.Item(3,41).Value = 11.5702479338843  # AO3
$ws.Cells.Item(3,42).Value = -0.2986111111111112  # AP3
$ws.Cells.Item(3,43).Value = 11.5702479338843  # AQ3

# --- Row 4 - Oman United Insurance Company SAOG (MSM:OUIC) ticker fix ---
$ws.Cells.Item(4,2).Value = "Oman United Insurance Company SAOG (MSM:OUIC)"  # B4
$ws.Cells.Item(4,4).Value = -0.0365  # D4
$ws.Cells.Item(4,5).Value = 0.174  # E4
$ws.Cells.Item(4,7).Value = 0.2319109461966605  # G4
$ws.Cells.Item(4,8).Value = 0.2319109461966605  # H4
$ws.Cells.Item(4,9).Value = 0.300556586270872  # I4
$ws.Cells.Item(4,10).Value = 0.25139146567718  # J4
$ws.Cells.Item(4,11).Value = 13.6  # K4
$ws.Cells.Item(4,12).Value = 0.2523191094619666  # L4
$ws.Cells.Item(4,13).Value = 7.79  # M4
$ws.Cells.Item(4,14).Value = 0.07497593840230991  # N4
$ws.Cells.Item(4,15).Value = 0.5727941176470588  # O4
$ws.Cells.Item(4,16).Value = 7.79  # P4
$ws.Cells.Item(4,17).Value = 0.07497593840230991  # Q4
$ws.Cells.Item(4,18).Value = 0.5727941176470588  # R4
$ws.Cells.Item(4,21).Value = 0  # U4
$ws.Cells.Item(4,22).Value = 0  # V4
$ws.Cells.Item(4,23).Value = 0.1912798874824191  # W4
$ws.Cells.Item(4,24).Value = 0.06370820384281146  # X4
$ws.Cells.Item(4,25).Value = 0.1275716836396077  # Y4
$ws.Cells.Item(4,26).Value = 0.7673143996013952  # Z4
$ws.Cells.Item(4,27).Value = 0.1928962915510001  # AA4
$ws.Cells.Item(4,28).Value = 0.06363327229597028  # AB4
$ws.Cells.Item(4,29).Value = 0.1292630192550298  # AC4
$ws.Cells.Item(4,30).Value = 0.384  # AD4
$ws.Cells.Item(4,32).Value = 0.384  # AF4
$ws.Cells.Item(4,33).Value = 0.384  # AG4
$ws.Cells.Item(4,34).Value = 0.003682252310996893  # AH4
$ws.Cells.Item(4,35).Value = 0.004962266101519694  # AI4
$ws.Cells.Item(4,36).Value = 0.003682252310996893  # AJ4
$ws.Cells.Item(4,37).Value = 0.004962266101519694  # AK4
$ws.Cells.Item(4,40).Value = 0.02232558139534884  # AN4
$ws.Cells.Item(4,42).Value = 0.02232558139534884  # AP4

# --- Row 5 - Al Madina Insurance Company SAOG (MSM:AMAT) ---
$ws.Cells.Item(5,4).Value = 0.0716  # D5
$ws.Cells.Item(5,5).Value = 0.278  # E5
$ws.Cells.Item(5,7).Value = 0.08249400479616306  # G5
$ws.Cells.Item(5,8).Value = 0.08249400479616306  # H5
$ws.Cells.Item(5,9).Value = 0.1021582733812949  # I5
$ws.Cells.Item(5,10).Value = 0.08419742345658357  # J5
$ws.Cells.Item(5,11).Value = 3.54  # K5
$ws.Cells.Item(5,12).Value = 0.08489208633093524  # L5
$ws.Cells.Item(5,13).Value = 2.73  # M5
$ws.Cells.Item(5,14).Value = 0.07690140845070423  # N5
$ws.Cells.Item(5,15).Value = 0.771186440677966  # O5
$ws.Cells.Item(5,16).Value = 2.73  # P5
$ws.Cells.Item(5,17).Value = 0.07690140845070423  # Q5
$ws.Cells.Item(5,18).Value = 0.771186440677966  # R5
$ws.Cells.Item(5,21).Value = 17.7  # U5
$ws.Cells.Item(5,22).Value = 0.4985915492957746  # V5
$ws.Cells.Item(5,23).Value = 0.05252225519287834  # W5
$ws.Cells.Item(5,24).Value = 0.06353781670736119  # X5
$ws.Cells.Item(5,25).Value = -0.01101556151448286  # Y5
$ws.Cells.Item(5,26).Value = 0.7433155080213903  # Z5
$ws.Cells.Item(5,27).Value = 0.06258525059072255  # AA5
$ws.Cells.Item(5,28).Value = 0.06353781670736119  # AB5
$ws.Cells.Item(5,29).Value = -0.0009525661166386451  # AC5
$ws.Cells.Item(5,33).Value = -17.7  # AG5
$ws.Cells.Item(5,36).Value = -0.99438202247191  # AJ5
$ws.Cells.Item(5,37).Value = -0.3582995951417005  # AK5
$ws.Cells.Item(5,42).Value = -3.75796178343949  # AP5

# --- Row 6 - NEW: Takaful Oman Insurance SAOG (MSM:TAOI) ---
$ws.Cells.Item(6,1).Value = "Oman"  # A6
$ws.Cells.Item(6,2).Value = "Takaful Oman Insurance SAOG (MSM:TAOI)"  # B6
$ws.Cells.Item(6,3).Value = "Insurance (General)"  # C6
$ws.Cells.Item(6,4).Value = 0.3779999999999999  # D6
$ws.Cells.Item(6,7).Value = -0.02977528089887641  # G6
$ws.Cells.Item(6,8).Value = -0.02977528089887641  # H6
$ws.Cells.Item(6,9).Value = -0.04466292134831461  # I6
$ws.Cells.Item(6,10).Value = -0.04466292134831461  # J6
$ws.Cells.Item(6,11).Value = -1.34  # K6
$ws.Cells.Item(6,12).Value = -0.03764044943820224  # L6
$ws.Cells.Item(6,13).Value = -0  # M6
$ws.Cells.Item(6,14).Value = -0  # N6
$ws.Cells.Item(6,15).Value = 0  # O6
$ws.Cells.Item(6,16).Value = -0  # P6
$ws.Cells.Item(6,17).Value = -0  # Q6
$ws.Cells.Item(6,18).Value = 0  # R6
$ws.Cells.Item(6,19).Value = 0  # S6
$ws.Cells.Item(6,21).Value = 0.353  # U6
$ws.Cells.Item(6,22).Value = 0.01124203821656051  # V6
$ws.Cells.Item(6,23).Value = -0.02567049808429119  # W6
$ws.Cells.Item(6,24).Value = 0.06830953545430658  # X6
$ws.Cells.Item(6,25).Value = -0.09398003353859777  # Y6
$ws.Cells.Item(6,26).Value = 0.6863973778077702  # Z6
$ws.Cells.Item(6,27).Value = -0.03065651209871783  # AA6
$ws.Cells.Item(6,28).Value = 0.06682872568239109  # AB6
$ws.Cells.Item(6,29).Value = -0.09748523778110892  # AC6
$ws.Cells.Item(6,30).Value = 3.25  # AD6
$ws.Cells.Item(6,31).Value = 0  # AE6
$ws.Cells.Item(6,32).Value = 3.25  # AF6
$ws.Cells.Item(6,33).Value = 2.897  # AG6
$ws.Cells.Item(6,34).Value = 0.09379509379509379  # AH6
$ws.Cells.Item(6,35).Value = 0.05220883534136546  # AI6
$ws.Cells.Item(6,36).Value = 0.08446802927369743  # AJ6
$ws.Cells.Item(6,37).Value = 0.0468035607541561  # AK6
$ws.Cells.Item(6,38).Value = 0  # AL6
$ws.Cells.Item(6,39).Value = 0  # AM6
$ws.Cells.Item(6,40).Value = -2.195945945945946  # AN6
$ws.Cells.Item(6,42).Value = -1.957432432432433  # AP6
